$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D price cells to be treated as Text so values like "1.00" or "198.60"
# are not silently reinterpreted as numbers (which would drop formatting/trailing zeros).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '76.307.03'
$ws.Range("E2").Value = '  +0.55%  '
$ws.Range("D3").Value = '2.964.58'
$ws.Range("E3").Value = '  +2.21%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '636.79'
$ws.Range("E5").Value = '  +6.68%  '
$ws.Range("D6").Value = '198.60'
$ws.Range("E6").Value = '  +0.20%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("E8").Value = '  -0.60%  '
$ws.Range("E9").Value = '  +2.33%  '
$ws.Range("D10").Value = '2.962.93'
$ws.Range("E10").Value = '  +2.21%  '
$ws.Range("D11").Value = '0.429'
$ws.Range("E11").Value = '  +2.34%  '
$ws.Range("E12").Value = '  -0.13%  '
$ws.Range("E13").Value = '  +2.05%  '
$ws.Range("D14").Value = '3.509.19'
$ws.Range("E14").Value = '  +2.53%  '
$ws.Range("D15").Value = '28.67'
$ws.Range("E15").Value = '  +4.60%  '
$ws.Range("D16").Value = '76.253.96'
$ws.Range("E16").Value = '  +0.53%  '
$ws.Range("E17").Value = '  -1.58%  '
$ws.Range("D18").Value = '2.960.35'
$ws.Range("E18").Value = '  +2.00%  '
$ws.Range("D19").Value = '13.31'
$ws.Range("E19").Value = '  +5.45%  '
$ws.Range("E20").Value = '  -2.13%  '
$ws.Range("D21").Value = '370.42'
$ws.Range("E21").Value = '  -1.63%  '
$ws.Range("E22").Value = '  +2.90%  '
$ws.Range("E23").Value = '  -3.21%  '
$ws.Range("D24").Value = '72.61'
$ws.Range("E24").Value = '  +2.23%  '
$ws.Range("D25").Value = '3.122.45'
$ws.Range("E25").Value = '  +2.43%  '
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  +0.19%  '
$ws.Range("D27").Value = '4.24'
$ws.Range("E27").Value = '  +0.76%  '
$ws.Range("D28").Value = '9.56'
$ws.Range("E28").Value = '  +0.35%  '
$ws.Range("E29").Value = '  -3.38%  '
$ws.Range("D30").Value = '0.999'
$ws.Range("E30").Value = '  +1.24%  '
$ws.Range("D31").Value = '8.22'
$ws.Range("E31").Value = '  +6.69%  '
$ws.Range("D32").Value = '512.35'
$ws.Range("E32").Value = '  +2.06%  '
$ws.Range("E33").Value = '  -0.93%  '
$ws.Range("D34").Value = '1.96'
$ws.Range("E34").Value = '  +9.13%  '
$ws.Range("E35").Value = '  +0.01%  '
$ws.Range("E36").Value = '  +0.10%  '
$ws.Range("E37").Value = '  +0.81%  '
$ws.Range("E38").Value = '  +1.34%  '
$ws.Range("D39").Value = '0.380'
$ws.Range("E39").Value = '  +11.56%  '
$ws.Range("D40").Value = '0.104'
$ws.Range("E40").Value = '  +15.59%  '
$ws.Range("E41").Value = '  -2.66%  '
$ws.Range("E42").Value = '  +0.03%  '
$ws.Range("D43").Value = '181.45'
$ws.Range("E43").Value = '  +0.88%  '
$ws.Range("D44").Value = '42.98'
$ws.Range("E44").Value = '  +7.43%  '
$ws.Range("D45").Value = '4.87'
$ws.Range("E45").Value = '  -2.50%  '
$ws.Range("D46").Value = '1.62'
$ws.Range("E46").Value = '  -2.16%  '
$ws.Range("E47").Value = '  +0.50%  '
$ws.Range("D48").Value = '0.702'
$ws.Range("E48").Value = '  +7.56%  '
$ws.Range("D49").Value = '0.580'
$ws.Range("E49").Value = '  +0.90%  '
$ws.Range("D50").Value = '2.27'
$ws.Range("E50").Value = '  -2.90%  '
$ws.Range("E51").Value = '  +2.04%  '
